# "biopython slides no longer have scripting"
#
# The trailing "Scripting Exercises" section (the title slide plus the
# Script #1 / Script #2 exercise slides) is removed from the deck.
# Those are the last 9 slides (positions 26-34) of the 34-slide deck;
# slides 1-25 are left completely untouched.

$p = $ppt.ActivePresentation

# Delete from the end so indices of not-yet-processed slides don't shift.
for ($i = $p.Slides.Count; $i -ge 26; $i--) {
    $p.Slides.Item($i).Delete()
}
